$d = $word.ActiveDocument

# The target paragraph currently reads (across 5 runs):
#   "Edit the " + "installChocoServer.ps1" + " script, if you wish to change the default " +
#   "$siteName" + ", and run it to install Chocolatey Server"
# It must become (across 2 runs):
#   "Once chocolatey is installed, run " + "choco install chocolatey.server"
#
# Locate each run's text range dynamically via Find (rather than hard-coded
# offsets) so the edit is robust to any earlier shifts in the document.

# Run 1 (plain formatting) - retext in place, keep its run properties.
$rEdit = $d.Content
$rEdit.Find.Execute("Edit the ") | Out-Null
$rEdit.Text = "Once chocolatey is installed, run "

# Run 2 (blue / underlined "installChocoServer.ps1") - remove entirely.
$rLink = $d.Content
$rLink.Find.Execute("installChocoServer.ps1") | Out-Null
$rLink.Text = ""

# Run 3 (" script, if you wish to change the default ") - remove entirely.
$rMid = $d.Content
$rMid.Find.Execute(" script, if you wish to change the default ") | Out-Null
$rMid.Text = ""

# Run 4 (Courier New "$siteName") - retext in place, keep its Courier New run properties.
$rSite = $d.Content
$rSite.Find.Execute("`$siteName") | Out-Null
$rSite.Text = "choco install chocolatey.server"

# Run 5 (", and run it to install Chocolatey Server") - remove entirely.
$rTail = $d.Content
$rTail.Find.Execute(", and run it to install Chocolatey Server") | Out-Null
$rTail.Text = ""
